$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the "Forandrad" (column C) date by one day for every data row (2-52)
for ($r = 2; $r -le 52; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value2 = $cCell.Value2 + 1
}

# Step 2: rows 7-52 got re-sorted/shuffled by the upstream data source.
# Re-assign Beteckning (A), Datum (B), Markagare (F) and Area (G) for each row
# to match the new order, row by row.

$ws.Cells.Item(7, 1).Value = "A 68624-2021"
$ws.Cells.Item(7, 2).Value = 44502
$ws.Cells.Item(7, 6).Value = $null
$ws.Cells.Item(7, 7).Value = 5.5

$ws.Cells.Item(8, 1).Value = "A 26702-2022"
$ws.Cells.Item(8, 2).Value = 44739.74783564815
$ws.Cells.Item(8, 6).Value = $null
$ws.Cells.Item(8, 7).Value = 3.8

$ws.Cells.Item(9, 1).Value = "A 29035-2021"
$ws.Cells.Item(9, 2).Value = 44358
$ws.Cells.Item(9, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item(9, 7).Value = 2.3

$ws.Cells.Item(11, 1).Value = "A 33953-2021"
$ws.Cells.Item(11, 2).Value = 44378
$ws.Cells.Item(11, 6).Value = $null
$ws.Cells.Item(11, 7).Value = 5

$ws.Cells.Item(13, 1).Value = "A 44069-2021"
$ws.Cells.Item(13, 2).Value = 44434
$ws.Cells.Item(13, 6).Value = $null
$ws.Cells.Item(13, 7).Value = 1.5

$ws.Cells.Item(14, 1).Value = "A 21264-2022"
$ws.Cells.Item(14, 2).Value = 44705
$ws.Cells.Item(14, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item(14, 7).Value = 2.4

$ws.Cells.Item(15, 1).Value = "A 21972-2023"
$ws.Cells.Item(15, 2).Value = 45068.66232638889
$ws.Cells.Item(15, 6).Value = $null
$ws.Cells.Item(15, 7).Value = 1.5

$ws.Cells.Item(16, 1).Value = "A 10263-2024"
$ws.Cells.Item(16, 2).Value = 45365.43090277778
$ws.Cells.Item(16, 6).Value = "Kyrkan"
$ws.Cells.Item(16, 7).Value = 1.4

$ws.Cells.Item(17, 1).Value = "A 30743-2021"
$ws.Cells.Item(17, 2).Value = 44365
$ws.Cells.Item(17, 6).Value = $null
$ws.Cells.Item(17, 7).Value = 3

$ws.Cells.Item(18, 1).Value = "A 55562-2022"
$ws.Cells.Item(18, 2).Value = 44888
$ws.Cells.Item(18, 6).Value = $null
$ws.Cells.Item(18, 7).Value = 0.8

$ws.Cells.Item(19, 1).Value = "A 27365-2025"
$ws.Cells.Item(19, 2).Value = 45812.64355324074
$ws.Cells.Item(19, 6).Value = $null
$ws.Cells.Item(19, 7).Value = 11.9

$ws.Cells.Item(20, 1).Value = "A 50230-2024"
$ws.Cells.Item(20, 2).Value = 45600
$ws.Cells.Item(20, 6).Value = $null
$ws.Cells.Item(20, 7).Value = 1.7

$ws.Cells.Item(21, 1).Value = "A 12077-2022"
$ws.Cells.Item(21, 2).Value = 44636.47484953704
$ws.Cells.Item(21, 6).Value = $null
$ws.Cells.Item(21, 7).Value = 2.1

$ws.Cells.Item(22, 1).Value = "A 37407-2023"
$ws.Cells.Item(22, 2).Value = 45156.60152777778
$ws.Cells.Item(22, 6).Value = $null
$ws.Cells.Item(22, 7).Value = 3.3

$ws.Cells.Item(23, 1).Value = "A 30174-2021"
$ws.Cells.Item(23, 2).Value = 44363
$ws.Cells.Item(23, 6).Value = $null
$ws.Cells.Item(23, 7).Value = 1.8

$ws.Cells.Item(24, 1).Value = "A 41546-2025"
$ws.Cells.Item(24, 2).Value = 45901.57927083333
$ws.Cells.Item(24, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item(24, 7).Value = 4.3

$ws.Cells.Item(25, 1).Value = "A 41550-2025"
$ws.Cells.Item(25, 2).Value = 45901.58652777778
$ws.Cells.Item(25, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item(25, 7).Value = 2.3

$ws.Cells.Item(27, 1).Value = "A 42991-2025"
$ws.Cells.Item(27, 2).Value = 45909.45190972222
$ws.Cells.Item(27, 6).Value = $null
$ws.Cells.Item(27, 7).Value = 7.2

$ws.Cells.Item(28, 1).Value = "A 42994-2025"
$ws.Cells.Item(28, 2).Value = 45909.45351851852
$ws.Cells.Item(28, 6).Value = $null
$ws.Cells.Item(28, 7).Value = 7.9

$ws.Cells.Item(29, 1).Value = "A 37072-2025"
$ws.Cells.Item(29, 2).Value = 45875.40431712963
$ws.Cells.Item(29, 6).Value = $null
$ws.Cells.Item(29, 7).Value = 1.5

$ws.Cells.Item(30, 1).Value = "A 37076-2025"
$ws.Cells.Item(30, 2).Value = 45875.41342592592
$ws.Cells.Item(30, 6).Value = $null
$ws.Cells.Item(30, 7).Value = 2.4

$ws.Cells.Item(31, 1).Value = "A 43448-2025"
$ws.Cells.Item(31, 2).Value = 45911.45209490741
$ws.Cells.Item(31, 6).Value = $null
$ws.Cells.Item(31, 7).Value = 1.2

$ws.Cells.Item(32, 1).Value = "A 13510-2025"
$ws.Cells.Item(32, 2).Value = 45736.47103009259
$ws.Cells.Item(32, 6).Value = $null
$ws.Cells.Item(32, 7).Value = 2.5

$ws.Cells.Item(33, 1).Value = "A 44192-2025"
$ws.Cells.Item(33, 2).Value = 45915.61556712963
$ws.Cells.Item(33, 6).Value = $null
$ws.Cells.Item(33, 7).Value = 0.8

$ws.Cells.Item(34, 1).Value = "A 14149-2022"
$ws.Cells.Item(34, 2).Value = 44651
$ws.Cells.Item(34, 6).Value = "Allmännings- och besparingsskogar"
$ws.Cells.Item(34, 7).Value = 3.8

$ws.Cells.Item(35, 1).Value = "A 53343-2024"
$ws.Cells.Item(35, 2).Value = 45614.43885416666
$ws.Cells.Item(35, 6).Value = $null
$ws.Cells.Item(35, 7).Value = 0.9

$ws.Cells.Item(36, 1).Value = "A 8848-2025"
$ws.Cells.Item(36, 2).Value = 45713.34292824074
$ws.Cells.Item(36, 6).Value = $null
$ws.Cells.Item(36, 7).Value = 2

$ws.Cells.Item(37, 1).Value = "A 8436-2023"
$ws.Cells.Item(37, 2).Value = 44977
$ws.Cells.Item(37, 6).Value = "Kyrkan"
$ws.Cells.Item(37, 7).Value = 4

$ws.Cells.Item(38, 1).Value = "A 7245-2025"
$ws.Cells.Item(38, 2).Value = 45702
$ws.Cells.Item(38, 6).Value = $null
$ws.Cells.Item(38, 7).Value = 4

$ws.Cells.Item(39, 1).Value = "A 34984-2024"
$ws.Cells.Item(39, 2).Value = 45527
$ws.Cells.Item(39, 6).Value = $null
$ws.Cells.Item(39, 7).Value = 4.1

$ws.Cells.Item(40, 1).Value = "A 58109-2025"
$ws.Cells.Item(40, 2).Value = 45982.59765046297
$ws.Cells.Item(40, 6).Value = $null
$ws.Cells.Item(40, 7).Value = 2.5

$ws.Cells.Item(41, 1).Value = "A 58111-2025"
$ws.Cells.Item(41, 2).Value = 45982.59920138889
$ws.Cells.Item(41, 6).Value = $null
$ws.Cells.Item(41, 7).Value = 0.6

$ws.Cells.Item(42, 1).Value = "A 58125-2025"
$ws.Cells.Item(42, 2).Value = 45982.61506944444
$ws.Cells.Item(42, 6).Value = $null
$ws.Cells.Item(42, 7).Value = 0.8

$ws.Cells.Item(43, 1).Value = "A 54207-2025"
$ws.Cells.Item(43, 2).Value = 45964
$ws.Cells.Item(43, 6).Value = $null
$ws.Cells.Item(43, 7).Value = 2.1

$ws.Cells.Item(44, 1).Value = "A 54203-2025"
$ws.Cells.Item(44, 2).Value = 45964
$ws.Cells.Item(44, 6).Value = $null
$ws.Cells.Item(44, 7).Value = 10.3

$ws.Cells.Item(45, 1).Value = "A 1621-2026"
$ws.Cells.Item(45, 2).Value = 46034.47645833333
$ws.Cells.Item(45, 6).Value = $null
$ws.Cells.Item(45, 7).Value = 1.4

$ws.Cells.Item(46, 1).Value = "A 1622-2026"
$ws.Cells.Item(46, 2).Value = 46034.47929398148
$ws.Cells.Item(46, 6).Value = $null
$ws.Cells.Item(46, 7).Value = 1.4

$ws.Cells.Item(47, 1).Value = "A 16762-2022"
$ws.Cells.Item(47, 2).Value = 44673.47876157407
$ws.Cells.Item(47, 6).Value = $null
$ws.Cells.Item(47, 7).Value = 4.2

$ws.Cells.Item(48, 1).Value = "A 22072-2023"
$ws.Cells.Item(48, 2).Value = 45069
$ws.Cells.Item(48, 6).Value = $null
$ws.Cells.Item(48, 7).Value = 3.5

$ws.Cells.Item(49, 1).Value = "A 62433-2025"
$ws.Cells.Item(49, 2).Value = 46007
$ws.Cells.Item(49, 6).Value = $null
$ws.Cells.Item(49, 7).Value = 2

$ws.Cells.Item(50, 1).Value = "A 30766-2022"
$ws.Cells.Item(50, 2).Value = 44764
$ws.Cells.Item(50, 6).Value = $null
$ws.Cells.Item(50, 7).Value = 0.6

$ws.Cells.Item(51, 1).Value = "A 35036-2024"
$ws.Cells.Item(51, 2).Value = 45527
$ws.Cells.Item(51, 6).Value = $null
$ws.Cells.Item(51, 7).Value = 1.7

$ws.Cells.Item(52, 1).Value = "A 62831-2023"
$ws.Cells.Item(52, 2).Value = 45270
$ws.Cells.Item(52, 6).Value = $null
$ws.Cells.Item(52, 7).Value = 1.1
